# Apply the GB_Transmission_Network_29_Bus.xlsx "gen" sheet change:
# Each grid interconnection (rows originally 67-71, buses 5, 10, 11, 26, 27)
# is split into two generator rows: one "import" generator (Pg_max=5000,
# Pg_min=0, positive cost coef H=50) and one "export" generator
# (Pg_max=0, Pg_min=-5000, zero cost coef H=0). Both rows are labeled with
# Gen Type "interconnector" in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gen")

# New layout: for each bus, one "import" row (max=5000,min=0,cost=50)
# followed by one "export" row (max=0,min=-5000,cost=0), starting at row 67.
$buses = @(5, 10, 11, 26, 27)

$targetRow = 67
foreach ($bus in $buses) {
    # Import generator row
    $ws.Cells.Item($targetRow, 1).Value = $bus
    $ws.Cells.Item($targetRow, 2).Value = "interconnector"
    $ws.Cells.Item($targetRow, 3).Value = 5000
    $ws.Cells.Item($targetRow, 4).Value = 0
    $ws.Cells.Item($targetRow, 5).Value = 0
    $ws.Cells.Item($targetRow, 6).Value = 0
    $ws.Cells.Item($targetRow, 7).Value = 0
    $ws.Cells.Item($targetRow, 8).Value = 50
    $targetRow++

    # Export generator row
    $ws.Cells.Item($targetRow, 1).Value = $bus
    $ws.Cells.Item($targetRow, 2).Value = "interconnector"
    $ws.Cells.Item($targetRow, 3).Value = 0
    $ws.Cells.Item($targetRow, 4).Value = -5000
    $ws.Cells.Item($targetRow, 5).Value = 0
    $ws.Cells.Item($targetRow, 6).Value = 0
    $ws.Cells.Item($targetRow, 7).Value = 0
    $ws.Cells.Item($targetRow, 8).Value = 0
    $targetRow++
}

# Update view state to mirror the saved selection/scroll position from the diff
$ws.Activate()
$ws.Range("C77").Select()
$excel.ActiveWindow.ScrollRow = 43
